$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common constant values for columns A,B,C,E,F,G,H,I,J across rows 283-304
$constA = 8
$constB = "Terminal La Palmera de La Serena"
$constC = "Coquimbo"
$constE = 4
$constF = "Fruta"
$constG = 100103
$constH = "Frutos de hueso (carozo)"
$constI = 100103004
$constJ = "Durazno"

$rows = @(
    @{ Row=283; D=44578; K='Carson'; L='Especial'; M=16; N=345000; O=350000; P=347500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=869; T=400 }
    @{ Row=284; D=44578; K='Carson'; L='Primera'; M=20; N=315000; O=320000; P=317500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=794; T=400 }
    @{ Row=285; D=44578; K='Carson'; L='Segunda'; M=20; N=275000; O=280000; P=277500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=694; T=400 }
    @{ Row=286; D=44266; K='Rich Lady'; L='Especial'; M=20; N=405000; O=410000; P=407500; Q='$/bins (400 kilos)'; R='Región Metropolitana'; S=1019; T=400 }
    @{ Row=287; D=44266; K='Rich Lady'; L='Primera'; M=20; N=375000; O=380000; P=377500; Q='$/bins (400 kilos)'; R='Región Metropolitana'; S=944; T=400 }
    @{ Row=288; D=44266; K='Kakamas'; L='Especial'; M=10; N=355000; O=360000; P=357500; Q='$/bins (400 kilos)'; R='Región Metropolitana'; S=894; T=400 }
    @{ Row=289; D=44266; K='Kakamas'; L='Primera'; M=20; N=315000; O=320000; P=317500; Q='$/bins (400 kilos)'; R='Región Metropolitana'; S=794; T=400 }
    @{ Row=290; D=44266; K='Kakamas'; L='Segunda'; M=16; N=275000; O=280000; P=277500; Q='$/bins (400 kilos)'; R='Región Metropolitana'; S=694; T=400 }
    @{ Row=291; D=44266; K='Phillips Cling'; L='Especial'; M=20; N=315000; O=320000; P=317500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=794; T=400 }
    @{ Row=292; D=44266; K='Phillips Cling'; L='Primera'; M=20; N=285000; O=290000; P=287500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=719; T=400 }
    @{ Row=293; D=44533; K='Florida King'; L='Primera'; M=16; N=405000; O=410000; P=407500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=1019; T=400 }
    @{ Row=294; D=44533; K='Florida King'; L='Segunda'; M=16; N=375000; O=380000; P=377500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=944; T=400 }
    @{ Row=295; D=44264; K='Phillips Cling'; L='Especial'; M=400; N=18500; O=19000; P=18750; Q='$/caja 16 kilos empedrada'; R='Región de O''Higgins'; S=1172; T=16 }
    @{ Row=296; D=44264; K='Phillips Cling'; L='Primera'; M=300; N=16500; O=17000; P=16750; Q='$/caja 16 kilos empedrada'; R='Región de O''Higgins'; S=1047; T=16 }
    @{ Row=297; D=44264; K='Phillips Cling'; L='Segunda'; M=200; N=13500; O=14000; P=13750; Q='$/caja 16 kilos empedrada'; R='Región de O''Higgins'; S=859; T=16 }
    @{ Row=298; D=44571; K='Carson'; L='Especial'; M=20; N=335000; O=340000; P=337500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=844; T=400 }
    @{ Row=299; D=44571; K='Carson'; L='Primera'; M=20; N=315000; O=320000; P=317500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=794; T=400 }
    @{ Row=300; D=44571; K='Toscana'; L='Primera'; M=16; N=325000; O=330000; P=327500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=819; T=400 }
    @{ Row=301; D=44279; K='Phillips Cling'; L='Primera'; M=22; N=325000; O=330000; P=327500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=819; T=400 }
    @{ Row=302; D=44279; K='Phillips Cling'; L='Segunda'; M=18; N=285000; O=290000; P=287500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=719; T=400 }
    @{ Row=303; D=44277; K='Phillips Cling'; L='Primera'; M=20; N=325000; O=330000; P=327500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=819; T=400 }
    @{ Row=304; D=44277; K='Phillips Cling'; L='Segunda'; M=20; N=285000; O=290000; P=287500; Q='$/bins (400 kilos)'; R='Región de O''Higgins'; S=719; T=400 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $constA
    $ws.Cells.Item($row, 2).Value = $constB
    $ws.Cells.Item($row, 3).Value = $constC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = $constE
    $ws.Cells.Item($row, 6).Value = $constF
    $ws.Cells.Item($row, 7).Value = $constG
    $ws.Cells.Item($row, 8).Value = $constH
    $ws.Cells.Item($row, 9).Value = $constI
    $ws.Cells.Item($row, 10).Value = $constJ
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}

Write-Output ("Final UsedRange=" + $ws.UsedRange.Address())
Write-Output ("Final dimension row count done")
